# Auto-generated edit script to apply the cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value as TEXT (avoids Excel auto-converting numeric-looking
# strings like '1.00' or '0.952' into numbers), while not leaving a permanent
# explicit cell style behind (keeps the default/general style afterwards).
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Plain text columns (B/C) never get mis-parsed as numbers, so assign directly.
function Set-PlainValue($cellRef, $value) {
    $ws.Range($cellRef).Value = $value
}

# Row 2
Set-TextValue "D2" '93.604.39'
Set-TextValue "E2" '  -4.24%  '

# Row 3
Set-TextValue "D3" '3.438.29'
Set-TextValue "E3" '  +2.12%  '

# Row 4
Set-TextValue "E4" '  -0.08%  '

# Row 5
Set-TextValue "D5" '236.36'
Set-TextValue "E5" '  -6.71%  '

# Row 6
Set-TextValue "D6" '639.52'
Set-TextValue "E6" '  -3.28%  '

# Row 7
Set-TextValue "E7" '  -0.74%  '

# Row 8
Set-TextValue "E8" '  -8.01%  '

# Row 9
Set-TextValue "E9" '  +0.05%  '

# Row 10
Set-TextValue "D10" '0.952'
Set-TextValue "E10" '  -7.02%  '

# Row 11
Set-TextValue "D11" '3.434.09'
Set-TextValue "E11" '  +1.98%  '

# Row 12
Set-TextValue "D12" '41.84'
Set-TextValue "E12" '  -0.36%  '

# Row 13
Set-TextValue "E13" '  -5.82%  '

# Row 14
Set-TextValue "D14" '6.14'
Set-TextValue "E14" '  +0.10%  '

# Row 15
Set-TextValue "D15" '4.085.85'
Set-TextValue "E15" '  +2.35%  '

# Row 16
Set-TextValue "D16" '93.511.05'
Set-TextValue "E16" '  -4.25%  '

# Row 17
Set-TextValue "D17" '0.0000249'
Set-TextValue "E17" '  -3.14%  '

# Row 18
Set-TextValue "E18" '  -5.83%  '

# Row 19
Set-TextValue "D19" '3.450.62'
Set-TextValue "E19" '  +2.40%  '

# Row 20
Set-TextValue "D20" '17.54'
Set-TextValue "E20" '  -2.69%  '

# Row 21
Set-TextValue "D21" '11.28'
Set-TextValue "E21" '  +3.58%  '

# Row 22
Set-TextValue "D22" '0.487'
Set-TextValue "E22" '  -9.81%  '

# Row 23
Set-TextValue "D23" '494.82'
Set-TextValue "E23" '  -3.80%  '

# Row 24
Set-TextValue "D24" '3.22'
Set-TextValue "E24" '  -4.99%  '

# Row 25
Set-TextValue "E25" '  -5.16%  '

# Row 26
Set-TextValue "E26" '  -5.85%  '

# Row 27
Set-TextValue "D27" '90.65'
Set-TextValue "E27" '  -6.59%  '

# Row 28
Set-TextValue "D28" '3.622.06'
Set-TextValue "E28" '  +2.16%  '

# Row 29
Set-TextValue "D29" '11.92'
Set-TextValue "E29" '  -4.33%  '

# Row 30
Set-TextValue "D30" '11.66'
Set-TextValue "E30" '  +0.08%  '

# Row 31
Set-TextValue "D31" '0.996'
Set-TextValue "E31" '  -0.66%  '

# Row 32
Set-TextValue "D32" '2.73'
Set-TextValue "E32" '  +5.02%  '

# Row 33
Set-TextValue "E33" '  -8.05%  '

# Row 34
Set-TextValue "E34" '  -6.91%  '

# Row 35
Set-TextValue "D35" '1.00'
Set-TextValue "E35" '  -0.05%  '

# Row 36
Set-TextValue "D36" '29.97'
Set-TextValue "E36" '  +3.77%  '

# Row 37
Set-TextValue "E37" '  -3.67%  '

# Row 38
Set-TextValue "D38" '550.29'
Set-TextValue "E38" '  +5.06%  '

# Row 39
Set-TextValue "D39" '7.56'
Set-TextValue "E39" '  -5.71%  '

# Row 40
Set-TextValue "E40" '  -5.07%  '

# Row 41
Set-TextValue "E41" '  -0.06%  '

# Row 42
Set-TextValue "E42" '  -1.38%  '

# Row 43
Set-TextValue "D43" '0.919'
Set-TextValue "E43" '  +6.49%  '

# Row 44
Set-TextValue "D44" '24.00'
Set-TextValue "E44" '  -1.78%  '

# Row 45
Set-TextValue "D45" '1.71'
Set-TextValue "E45" '  -2.06%  '

# Row 46
Set-TextValue "D46" '0.0410'
Set-TextValue "E46" '  -7.00%  '

# Row 47
Set-PlainValue "B47" 'Filecoin'
Set-PlainValue "C47" 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue "D47" '5.50'
Set-TextValue "E47" '  -3.79%  '

# Row 48
Set-PlainValue "B48" 'MantraDAO'
Set-PlainValue "C48" 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
Set-TextValue "D48" '3.53'
Set-TextValue "E48" '  -3.08%  '

# Row 49
Set-TextValue "E49" '  +4.16%  '

# Row 50
Set-TextValue "E50" '  +0.59%  '

# Row 51
Set-TextValue "D51" '52.96'
Set-TextValue "E51" '  -1.01%  '
